# Update "想去人数" (interest count, column F) figures to their latest
# scraped values across the three sheets that contain event data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 37
$ws.Range("F3").Value = 505
$ws.Range("F5").Value = 101
$ws.Range("F6").Value = 342
$ws.Range("F7").Value = 1321
$ws.Range("F8").Value = 526
$ws.Range("F10").Value = 1366
$ws.Range("F11").Value = 187
$ws.Range("F13").Value = 205
$ws.Range("F14").Value = 137
$ws.Range("F15").Value = 269
$ws.Range("F16").Value = 1713
$ws.Range("F17").Value = 634
$ws.Range("F18").Value = 282
$ws.Range("F19").Value = 320
$ws.Range("F20").Value = 3512
$ws.Range("F21").Value = 31
$ws.Range("F23").Value = 945
$ws.Range("F24").Value = 1226
$ws.Range("F26").Value = 2871
$ws.Range("F27").Value = 1694
$ws.Range("F29").Value = 25
$ws.Range("F30").Value = 135
$ws.Range("F32").Value = 882
$ws.Range("F33").Value = 13
$ws.Range("F34").Value = 2000
$ws.Range("F36").Value = 2024
$ws.Range("F38").Value = 488
$ws.Range("F39").Value = 189
$ws.Range("F41").Value = 355
$ws.Range("F42").Value = 945
$ws.Range("F43").Value = 825
$ws.Range("F44").Value = 1078
$ws.Range("F45").Value = 182
$ws.Range("F46").Value = 459
$ws.Range("F47").Value = 302
$ws.Range("F48").Value = 243
$ws.Range("F49").Value = 3393
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 2
$ws.Range("F13").Value = 828
$ws.Range("F22").Value = 19
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 505
$ws.Range("F4").Value = 101
$ws.Range("F7").Value = 1321
$ws.Range("F8").Value = 1366
$ws.Range("F9").Value = 187
$ws.Range("F11").Value = 205
$ws.Range("F12").Value = 137
$ws.Range("F14").Value = 269
$ws.Range("F15").Value = 1713
$ws.Range("F16").Value = 634
$ws.Range("F17").Value = 282
$ws.Range("F18").Value = 320
$ws.Range("F19").Value = 3513
$ws.Range("F20").Value = 31
$ws.Range("F23").Value = 1226
$ws.Range("F24").Value = 2871
$ws.Range("F25").Value = 1694
$ws.Range("F27").Value = 135
$ws.Range("F28").Value = 828
$ws.Range("F31").Value = 882
$ws.Range("F32").Value = 2000
$ws.Range("F35").Value = 2024
$ws.Range("F36").Value = 488
$ws.Range("F37").Value = 189
$ws.Range("F39").Value = 945
$ws.Range("F40").Value = 825
$ws.Range("F41").Value = 1078
$ws.Range("F42").Value = 182
$ws.Range("F43").Value = 459
$ws.Range("F44").Value = 302
$ws.Range("F46").Value = 19
$ws.Range("F47").Value = 243
$ws.Range("F48").Value = 3393
